$d = $word.ActiveDocument

# 1. Title: "... Healthy Older Adults from the Knight ADRC" -> "... Healthy Older
#    Adults from OASIS and ADNI"
$d.Content.Find.Execute("the Knight ADRC", $true, $false, $false, $false, $false, `
    $true, 1, $false, "OASIS and ADNI", 2)

# 2. Body paragraph: swap the Knight ADRC cohort description for the OASIS one,
#    keeping the trailing "and the Alzheimer Disease Neuroimaging Initiative
#    (ADNI) cohorts" text untouched.
$d.Content.Find.Execute( `
    "Charles F. and Joanne Knight Alzheimer Disease Research Center (Knight ADRC)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Open Access Series of Imaging Studies (OASIS)", 2)
